# "update project to use tree-based cell type mapping"
#
# The "cell_types" sheet maps a cell-type label (column A) to a color
# (column B). This renames the cell-type labels to the new, tree-based
# naming scheme (colors / row order / styles are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cell_types")

$ws.Range("A2").Value2  = "T cell CD4+ (non-regulatory)"
$ws.Range("A3").Value2  = "T cell CD8+"
$ws.Range("A4").Value2  = "NK cell"
$ws.Range("A5").Value2  = "Dendritic cell"
$ws.Range("A6").Value2  = "T cell regulatory (Tregs)"
$ws.Range("A7").Value2  = "Macrophage/Monocyte"
$ws.Range("A8").Value2  = "B cell"
$ws.Range("A9").Value2  = "Endothelial cell"
$ws.Range("A10").Value2 = "Cancer associated fibroblast"
$ws.Range("A11").Value2 = "Melanoma cell"
$ws.Range("A12").Value2 = "Ovarian carcinoma cell"
$ws.Range("A13").Value2 = "other cell"

# Move/record the active selection on the cell_types sheet, as in the diff.
$ws.Activate()
$ws.Range("A3").Select()
